# Update cryptocurrency price/volume data to the values scraped on 2023-01-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (cells hold text, not numbers, so force a text
# number format before assigning to avoid Excel auto-converting "300.88" to a
# number or "0.52%" to a percentage value).
$updates = @{
    "D2" = "300.88"
    "E2" = "0.52%"
    "D3" = "31.50"
    "E3" = "1.12%"
    "D4" = "5.084"
    "E4" = "-0.97%"
    "D5" = "0.07848"
    "E5" = "-3.59%"
    "D6" = "2.315"
    "E6" = "-5.18%"
    "D7" = "7.809"
    "E7" = "-0.39%"
    "D8" = "3.834"
    "E8" = "-0.23%"
    "D9" = "0.9176"
    "E9" = "0.65%"
    "D10" = "0.1752"
    "E10" = "2.50%"
    "D11" = "0.07575"
    "E11" = "3.65%"
    "D12" = "0.09279"
    "E12" = "15.73%"
    "D13" = "0.02993"
    "E13" = "-1.10%"
    "E14" = "0.60%"
    "D15" = "0.001509"
    "E15" = "0.29%"
    "D16" = "0.005788"
    "E16" = "-3.24%"
    "D17" = "3.469"
    "E17" = "-0.77%"
    "D18" = "2.248"
    "E19" = "-1.13%"
    "D20" = "0.1327"
    "E20" = "-1.50%"
    "D21" = "4.064"
    "E21" = "-11.90%"
    "D22" = "0.1790"
    "E22" = "11.59%"
    "E23" = "0.51%"
    "D24" = "0.001252"
    "E24" = "-0.73%"
    "D25" = "0.004464"
    "E25" = "0.42%"
    "D26" = "0.0001250"
    "E26" = "5.62%"
    "E27" = "-1.63%"
    "D39" = "0.01757"
    "E39" = "-3.06%"
    "D40" = "0.04704"
    "E40" = "3.79%"
    "D41" = "0.007163"
    "E41" = "-0.84%"
    "D42" = "0.1361"
    "E42" = "1.40%"
    "D43" = "0.002189"
    "E43" = "0.13%"
    "D44" = "0.009777"
    "E44" = "-8.71%"
    "D45" = "0.00006270"
    "E45" = "0.14%"
    "D46" = "0.00000000750"
    "E46" = "-0.27%"
    "E47" = "19.96%"
    "D48" = "0.7435"
    "E48" = "-9.39%"
    "D49" = "0.00002099"
    "E49" = "-0.27%"
    "D50" = "0.0001999"
    "E50" = "-0.27%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
